# Generate Report for Archive
#
# 1) Change the "Ready for handoff" status text to "In Translation" everywhere
#    it appears (Overview sheet zh-cn/de-de columns, and the Status column on
#    the zh-cn and de-de sheets).
# 2) Narrow the "Status" column width (Overview!E:F, zh-cn!C, de-de!C) from
#    ~17.22 characters down to ~13.41 characters.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells, rows 2-3 ---
foreach ($row in 2..3) {
    foreach ($col in @(5, 6)) {
        $cell = $wsOverview.Cells.Item($row, $col)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}

# --- zh-cn / de-de sheets: Status column (col C), rows 2-3 ---
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($row in 2..3) {
        $cell = $ws.Cells.Item($row, 3)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}

# --- Column width changes ---
# Target stored (OOXML) column width is 13.4101848602295 characters. The
# engine quantizes ColumnWidth to 1/6-character increments when it persists
# it, so request the value whose quantized result lands closest to the
# target (12.5 -> stored width 13.333333333333334, the nearest reachable
# value to 13.4101848602295).
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
